$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/16/2025  Through  6/22/2025"

# --- Data table updates (rows 14-30) ---
# Row 14
$ws.Range("F14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 2
$ws.Range("I14").Value = 10
$ws.Range("K14").Value = 233.333333333333
$ws.Range("N14").Value = -47.368421052631

# Row 15
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("G15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("I15").Value = 15
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = 15.384615384615
$ws.Range("L15").Value = -11.764705882352
$ws.Range("M15").Value = 36.363636363636
$ws.Range("N15").Value = -50

# Row 16
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -36.363636363636
$ws.Range("F16").Value = 45
$ws.Range("H16").Value = 28.571428571428
$ws.Range("I16").Value = 213
$ws.Range("J16").Value = 231
$ws.Range("K16").Value = -7.792207792207
$ws.Range("L16").Value = 23.121387283237
$ws.Range("M16").Value = 29.090909090909
$ws.Range("N16").Value = -61.131386861313

# Row 17
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 28
$ws.Range("E17").Value = -53.571428571428
$ws.Range("G17").Value = 80
$ws.Range("H17").Value = -30
$ws.Range("I17").Value = 354
$ws.Range("J17").Value = 355
$ws.Range("K17").Value = -0.281690140845
$ws.Range("L17").Value = 14.193548387096
$ws.Range("M17").Value = 118.518518518519
$ws.Range("N17").Value = 10.625

# Row 18
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -22.222222222222
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 31
$ws.Range("H18").Value = -12.903225806451
$ws.Range("I18").Value = 133
$ws.Range("J18").Value = 149
$ws.Range("K18").Value = -10.738255033557
$ws.Range("L18").Value = 27.884615384615
$ws.Range("M18").Value = 22.018348623853
$ws.Range("N18").Value = -79.153605015674

# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 18.181818181818
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 11.904761904761
$ws.Range("I19").Value = 290
$ws.Range("J19").Value = 292
$ws.Range("K19").Value = -0.684931506849
$ws.Range("L19").Value = 28.888888888888
$ws.Range("M19").Value = 97.278911564625
$ws.Range("N19").Value = 42.156862745098

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = -58.333333333333
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 127
$ws.Range("J20").Value = 137
$ws.Range("K20").Value = -7.299270072992
$ws.Range("L20").Value = -19.108280254777
$ws.Range("M20").Value = 170.212765957447
$ws.Range("N20").Value = -54.804270462633

# Row 21
$ws.Range("C21").Value = 47
$ws.Range("D21").Value = 72
$ws.Range("E21").Value = -34.722222222222
$ws.Range("F21").Value = 206
$ws.Range("G21").Value = 220
$ws.Range("H21").Value = -6.363636363636
$ws.Range("I21").Value = 1142
$ws.Range("J21").Value = 1180
$ws.Range("K21").Value = -3.22033898305
$ws.Range("L21").Value = 14.658634538152
$ws.Range("M21").Value = 78.159126365054
$ws.Range("N21").Value = -44.019607843137

# Row 22
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("M22").Value = 25

# Row 23
$ws.Range("D23").Value = 2
$ws.Range("J23").Value = 15
$ws.Range("K23").Value = -66.666666666666
$ws.Range("L23").Value = -54.545454545454

# Row 24
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = 146.153846153846
$ws.Range("F24").Value = 112
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = 86.666666666666
$ws.Range("I24").Value = 630
$ws.Range("J24").Value = 477
$ws.Range("K24").Value = 32.075471698113
$ws.Range("L24").Value = 27.016129032258
$ws.Range("M24").Value = 43.181818181818

# Row 25
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 55.555555555555
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 66.666666666666
$ws.Range("I25").Value = 155
$ws.Range("J25").Value = 121
$ws.Range("K25").Value = 28.099173553719
$ws.Range("L25").Value = -9.356725146198

# Row 26
$ws.Range("C26").Value = 33
$ws.Range("D26").Value = 23
$ws.Range("E26").Value = 43.478260869565
$ws.Range("F26").Value = 108
$ws.Range("G26").Value = 102
$ws.Range("H26").Value = 5.882352941176
$ws.Range("I26").Value = 522
$ws.Range("J26").Value = 463
$ws.Range("K26").Value = 12.742980561555
$ws.Range("L26").Value = 30.827067669172
$ws.Range("M26").Value = 12.987012987013

# Row 27
$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -40
$ws.Range("I27").Value = 22
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = -8.333333333333
$ws.Range("L27").Value = -24.137931034482

# Row 28
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("A28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("A28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("A28").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 20
$ws.Range("I28").Value = 41
$ws.Range("K28").Value = 24.242424242424
$ws.Range("L28").Value = 28.125

# Row 29
$ws.Range("D29").Value = 1
$ws.Range("G29").Value = 8
$ws.Range("J29").Value = 18
$ws.Range("K29").Value = -33.333333333333
$ws.Range("L29").Value = -36.842105263157
$ws.Range("M29").Value = -36.842105263157
$ws.Range("N29").Value = -76.923076923076

# Row 30
$ws.Range("D30").Value = 1
$ws.Range("G30").Value = 6
$ws.Range("J30").Value = 16
$ws.Range("K30").Value = -25
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -25
$ws.Range("N30").Value = -72.093023255813
